$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
